# Updates the cryptos list: refreshed prices/1h-volume figures for all
# existing coins, plus a reshuffle of the bottom of the table where
# BabyDogeCoin enters the rankings (displacing Aptos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as TEXT, matching the sheet's existing inline-string
# cells. Values such as "0.9981" or "26.072.50" would otherwise be auto-
# detected by Excel as numbers, so numeric-looking strings are written with a
# leading apostrophe (forces text) and the style is reset to Normal right
# after so no stray "quote prefix" cell format is left behind.
function Set-TextValue {
    param($cell, [string]$val)
    $isNumericLooking = $val -match '^[+-]?(\d+\.?\d*|\.\d+)([eE][+-]?\d+)?$'
    if ($isNumericLooking) {
        $cell.Value = "'" + $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

# --- Rows 2-45: refreshed Price (D) / Volume(1h) (E) values ---
Set-TextValue $ws.Range("D2") "26.072.50"
Set-TextValue $ws.Range("E2") "  -0.28%  "
Set-TextValue $ws.Range("D3") "1.645.55"
Set-TextValue $ws.Range("E3") "  -1.41%  "
Set-TextValue $ws.Range("E4") "  -0.20%  "
Set-TextValue $ws.Range("D5") "215.52"
Set-TextValue $ws.Range("E5") "  +2.35%  "
Set-TextValue $ws.Range("D6") "0.5219"
Set-TextValue $ws.Range("E6") "  +0.42%  "
Set-TextValue $ws.Range("E7") "  -0.16%  "
Set-TextValue $ws.Range("D8") "0.2612"
Set-TextValue $ws.Range("E8") "  -0.43%  "
Set-TextValue $ws.Range("D9") "0.06374"
Set-TextValue $ws.Range("E9") "  +0.83%  "
Set-TextValue $ws.Range("D10") "20.83"
Set-TextValue $ws.Range("E10") "  -1.59%  "
Set-TextValue $ws.Range("D11") "0.07665"
Set-TextValue $ws.Range("E11") "  +1.65%  "
Set-TextValue $ws.Range("D12") "1.648.73"
Set-TextValue $ws.Range("E12") "  -1.30%  "
Set-TextValue $ws.Range("D13") "4.422"
Set-TextValue $ws.Range("E13") "  -0.45%  "
Set-TextValue $ws.Range("D14") "1.868.54"
Set-TextValue $ws.Range("E14") "  -1.57%  "
Set-TextValue $ws.Range("D15") "0.5548"
Set-TextValue $ws.Range("E15") "  +1.10%  "
Set-TextValue $ws.Range("D16") "0.0₅8257"
Set-TextValue $ws.Range("E16") "  +3.19%  "
Set-TextValue $ws.Range("D17") "64.99"
Set-TextValue $ws.Range("E17") "  -2.01%  "
Set-TextValue $ws.Range("D18") "26.088.33"
Set-TextValue $ws.Range("E19") "  -0.12%  "
Set-TextValue $ws.Range("D20") "4.729"
Set-TextValue $ws.Range("E20") "  -0.53%  "
Set-TextValue $ws.Range("D21") "188.42"
Set-TextValue $ws.Range("E21") "  +0.76%  "
Set-TextValue $ws.Range("D22") "10.21"
Set-TextValue $ws.Range("E22") "  -1.10%  "
Set-TextValue $ws.Range("D23") "6.238"
Set-TextValue $ws.Range("E23") "  +0.44%  "
Set-TextValue $ws.Range("E24") "  -0.20%  "
Set-TextValue $ws.Range("D25") "146.17"
Set-TextValue $ws.Range("E25") "  -2.56%  "
Set-TextValue $ws.Range("D26") "0.1221"
Set-TextValue $ws.Range("E26") "  -1.27%  "
Set-TextValue $ws.Range("D27") "7.423"
Set-TextValue $ws.Range("E27") "  -0.81%  "
Set-TextValue $ws.Range("D28") "15.85"
Set-TextValue $ws.Range("E28") "  +0.22%  "
Set-TextValue $ws.Range("D29") "1.383"
Set-TextValue $ws.Range("E29") "  +2.54%  "
Set-TextValue $ws.Range("D30") "0.05955"
Set-TextValue $ws.Range("E30") "  -5.46%  "
Set-TextValue $ws.Range("E31") "  -1.34%  "
Set-TextValue $ws.Range("D32") "3.399"
Set-TextValue $ws.Range("E32") "  -0.32%  "
Set-TextValue $ws.Range("D33") "3.404"
Set-TextValue $ws.Range("E33") "  -3.26%  "
Set-TextValue $ws.Range("E34") "  +1.26%  "
Set-TextValue $ws.Range("D35") "0.9937"
Set-TextValue $ws.Range("E35") "  -1.03%  "
Set-TextValue $ws.Range("E36") "  -0.41%  "
Set-TextValue $ws.Range("D37") "2.751"
Set-TextValue $ws.Range("E37") "  -0.46%  "
Set-TextValue $ws.Range("D38") "0.5634"
Set-TextValue $ws.Range("E38") "  -6.95%  "
Set-TextValue $ws.Range("D39") "0.01617"
Set-TextValue $ws.Range("E39") "  +0.26%  "
Set-TextValue $ws.Range("D40") "0.8586"
Set-TextValue $ws.Range("E40") "  -0.76%  "
Set-TextValue $ws.Range("D41") "5.836"
Set-TextValue $ws.Range("E41") "  -4.45%  "
Set-TextValue $ws.Range("E42") "  -0.26%  "
Set-TextValue $ws.Range("D43") "1.029.40"
Set-TextValue $ws.Range("E43") "  -7.38%  "
Set-TextValue $ws.Range("D44") "99.31"
Set-TextValue $ws.Range("E44") "  -1.16%  "
Set-TextValue $ws.Range("D45") "1.794.98"
Set-TextValue $ws.Range("E45") "  -1.51%  "

# --- Rows 46-51: BabyDogeCoin enters the table, Aptos drops out, and the
#     remaining coins (Aave, Frax, EnergySwap, Cronos, Mantle) shift / refresh ---
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D46") "0.0₈112"
Set-TextValue $ws.Range("E46") "  +1.45%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "55.85"
Set-TextValue $ws.Range("E47") "  +0.67%  "
Set-TextValue $ws.Range("D48") "0.9981"
Set-TextValue $ws.Range("E48") "  +0.08%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "8.073"
Set-TextValue $ws.Range("E49") "  +0.15%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.05158"
Set-TextValue $ws.Range("E50") "  -1.54%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D51") "0.4222"
Set-TextValue $ws.Range("E51") "  -0.49%  "
